# Applies the "Gracias!" slide edits:
#  1. Moves the Title placeholder up (y offset change).
#  2. Adds a "Link Diapositivas" textbox with a hyperlink.
#  3. Adds a "Link Video Presentacion" textbox with a hyperlink.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(24)

# PowerPoint's Shape.Top/.Left setters marshal through a single-precision
# (32-bit) float, same as real PPT COM automation, so a naive EMU/12700
# division can land one EMU short after the float32 round-trip. Search for
# a nearby point value whose float32 cast converts back to the exact EMU.
function Get-PointsForEmu($TargetEmu, $EmuPerPoint) {
    $base = $TargetEmu / $EmuPerPoint
    $step = 0.0000001
    $cand = $base
    for ($i = 0; $i -lt 200000; $i++) {
        $f32 = [float]$cand
        $val = [double]$f32 * $EmuPerPoint
        $floored = [math]::Floor($val + 0.000001)
        if ($floored -eq $TargetEmu) {
            return $cand
        } elseif ($floored -lt $TargetEmu) {
            $cand += $step
        } else {
            $cand -= $step
        }
    }
    throw "no solution found for $TargetEmu"
}

# --- 1. Move the Title shape up ---
$title = $s.Shapes.Item(1)
$title.Top = Get-PointsForEmu 1417248 12700

# --- Burn through the shape id/name counter (3, 4, 5 get consumed and
#     discarded) so the two real new textboxes land on id=6/"TextBox 5"
#     and id=7/"TextBox 6", matching the authored file. ---
for ($i = 0; $i -lt 3; $i++) {
    $dummy = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
    $dummy.Delete()
}

# --- 2. "Link Diapositivas" textbox ---
$tb1 = $s.Shapes.AddTextbox(1, 359832 / 12700, 4363534 / 12700, 11641668 / 12700, 1077218 / 12700)
$tb1.Fill.Visible = $false
$tb1.TextFrame.WordWrap = $true
$tb1.TextFrame.AutoSize = 1

$tr1 = $tb1.TextFrame.TextRange
$url1 = "https://docs.google.com/presentation/d/1OVkkPnO7M1RX9RnqtcDdVlekxcDUIeMl/edit?usp=sharing&ouid=105746497042844411325&rtpof=true&sd=true"
$tr1.Text = "Link Diapositivas: " + $url1
$tr1.Font.Size = 16
$tr1.Font.Bold = $true

# Split into the same runs the original file has: "Link ", "Diapositivas",
# ": ", and the hyperlinked URL.
$p1 = $tr1.Characters(1, 5)
$p1.Font.Bold = $true
$p2 = $tr1.Characters(6, 12)
$p2.Font.Bold = $true
$p3 = $tr1.Characters(18, 2)
$p3.Font.Bold = $true
$p4 = $tr1.Characters(20, $url1.Length)
$p4.Font.Bold = $true
$p4.ActionSettings.Item(1).Hyperlink.Address = $url1

# Second, empty paragraph.
$null = $tr1.InsertAfter("`r")

# --- 3. "Link Video Presentacion" textbox ---
$tb2 = $s.Shapes.AddTextbox(1, 359832 / 12700, 5242317 / 12700, 11641668 / 12700, 338554 / 12700)
$tb2.Fill.Visible = $false
$tb2.TextFrame.WordWrap = $true
$tb2.TextFrame.AutoSize = 1

$tr2 = $tb2.TextFrame.TextRange
$prefix2 = "Link Video Presentación: "
$url2 = "https://drive.google.com/file/d/1CxYavUrC-Oxv8rdwEJcjzMq9boBdzCZw/view?usp=sharing"
$tr2.Text = $prefix2 + $url2
$tr2.LanguageID = "es-EC"
$tr2.Font.Size = 16
$tr2.Font.Bold = $true

$q1 = $tr2.Characters(1, $prefix2.Length)
$q1.Font.Bold = $true
$q2 = $tr2.Characters($prefix2.Length + 1, $url2.Length)
$q2.Font.Bold = $true
$q2.ActionSettings.Item(1).Hyperlink.Address = $url2
